# Apply the 20260204 monitor-xlsx update:
#  - overview ("總覽") + details ("詳細數據") sheets: refreshed market/chips snapshot text values
#  - stock-chips ("個股籌碼") sheet: refreshed numeric columns, plus two newly tracked metrics
#    (K = 融資增減(張) margin-financing delta, M = 借券增減(張) securities-lending delta)

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("總覽")
$wsDetail   = $wb.Worksheets.Item("詳細數據")
$wsChips    = $wb.Worksheets.Item("個股籌碼")

# --- 總覽 (Overview): text snapshot cells -------------------------------------
# Values that Excel would otherwise auto-parse as numbers/percentages need a
# leading quote so they stay literal text, matching the sheet's existing style.

$wsOverview.Range("C4").Value = "'4.28%"
$wsOverview.Range("D4").Value = "'+0.14%"
$wsOverview.Range("C5").Value = "'5034.4`$"
$wsOverview.Range("D5").Value = "'+2.67%"
$wsOverview.Range("C6").Value = "'31.57"
$wsOverview.Range("D6").Value = "'-0.02%"
$wsOverview.Range("C7").Value = "'7619.16"
$wsOverview.Range("D7").Value = "'-4.36%"
$wsOverview.Range("C8").Value = "'18.64"
$wsOverview.Range("D8").Value = "'+3.56%"
$wsOverview.Range("C10").Value = "58.21億"
$wsOverview.Range("E10").Value = "30.26億"
$wsOverview.Range("F10").Value = "151.3億"
$wsOverview.Range("C12").Value = "162.17億"
$wsOverview.Range("C14").Value = "-24.38億"
$wsOverview.Range("E14").Value = "-4.08億"
$wsOverview.Range("F14").Value = "-20.38億"
$wsOverview.Range("G14").Value = "-18.19億"
$wsOverview.Range("H14").Value = "-363.83億"

# --- 詳細數據 (Details): mirrors the same 5D/20D aggregates shown on 總覽 -----

$wsDetail.Range("B9").Value = "30.26億"
$wsDetail.Range("B10").Value = "151.3億"
$wsDetail.Range("B14").Value = "-4.08億"
$wsDetail.Range("B15").Value = "-20.38億"
$wsDetail.Range("B16").Value = "-18.19億"
$wsDetail.Range("B17").Value = "-363.83億"

# --- 個股籌碼 (Stock chips): per-ticker numeric refresh ------------------------

# row 4
$wsChips.Range("K4").Value = 7630
$wsChips.Range("M4").Value = -4012076
$wsChips.Range("R4").Value = 137907000
$wsChips.Range("S4").Value = 0.99

# row 5
$wsChips.Range("K5").Value = 10516
$wsChips.Range("M5").Value = -29673
$wsChips.Range("R5").Value = 2003000

# row 6
$wsChips.Range("K6").Value = 2821
$wsChips.Range("M6").Value = -78691
$wsChips.Range("R6").Value = 3409000
$wsChips.Range("S6").Value = 0.8

# row 7
$wsChips.Range("E7").Value = 82828
$wsChips.Range("K7").Value = 83558
$wsChips.Range("M7").Value = -1105334
$wsChips.Range("R7").Value = 176374918
$wsChips.Range("S7").Value = 0.68

# row 8
$wsChips.Range("E8").Value = 7342
$wsChips.Range("K8").Value = 6152
$wsChips.Range("M8").Value = -648628
$wsChips.Range("R8").Value = 2200680
$wsChips.Range("S8").Value = 0.18

# row 9
$wsChips.Range("E9").Value = 29087
$wsChips.Range("H9").Value = 840
$wsChips.Range("K9").Value = 22274
$wsChips.Range("M9").Value = -6482877
$wsChips.Range("R9").Value = 1556680
$wsChips.Range("S9").Value = 0.04
$wsChips.Range("T9").Value = 1746.76

# row 10
$wsChips.Range("E10").Value = 210454
$wsChips.Range("K10").Value = 143776
$wsChips.Range("M10").Value = -1110204
$wsChips.Range("R10").Value = 31176000
$wsChips.Range("S10").Value = 0.11

# row 11
$wsChips.Range("E11").Value = 6355
$wsChips.Range("K11").Value = 2160
$wsChips.Range("M11").Value = -89471
$wsChips.Range("R11").Value = 1811528
$wsChips.Range("S11").Value = 0.45
$wsChips.Range("T11").Value = 1714.99
$wsChips.Range("U11").Value = 18.95

# row 12
$wsChips.Range("E12").Value = 2252
$wsChips.Range("K12").Value = 5915
$wsChips.Range("M12").Value = -19723
$wsChips.Range("R12").Value = 5442503
$wsChips.Range("S12").Value = 2.38

# row 13
$wsChips.Range("K13").Value = 28536
$wsChips.Range("M13").Value = -267204
$wsChips.Range("R13").Value = 22195945
$wsChips.Range("S13").Value = 0.44

# row 14
$wsChips.Range("H14").Value = 17
$wsChips.Range("K14").Value = 4429
$wsChips.Range("M14").Value = -19211
$wsChips.Range("R14").Value = 1682000
$wsChips.Range("S14").Value = 0.33

# row 15
$wsChips.Range("K15").Value = -50
$wsChips.Range("M15").Value = 0
$wsChips.Range("R15").Value = 2188400
$wsChips.Range("S15").Value = 3.13

# row 16
$wsChips.Range("K16").Value = 2454
$wsChips.Range("M16").Value = 0
$wsChips.Range("R16").Value = 8566640
$wsChips.Range("S16").Value = 0.45

# row 17
$wsChips.Range("K17").Value = 22
$wsChips.Range("M17").Value = 0
$wsChips.Range("R17").Value = 1887000
$wsChips.Range("S17").Value = 1.94

# row 18
$wsChips.Range("K18").Value = 198
$wsChips.Range("M18").Value = 0
$wsChips.Range("R18").Value = 11349993
$wsChips.Range("S18").Value = 0.55

# row 19
$wsChips.Range("K19").Value = 25655
$wsChips.Range("M19").Value = -112769
$wsChips.Range("R19").Value = 10599826
$wsChips.Range("S19").Value = 0.18

